# Apply updated cryptocurrency market data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.328.44"
$ws.Range("E2").Value = "  -2.85%  "
# Row 3
$ws.Range("D3").Value = "1.940.50"
$ws.Range("E3").Value = "  -3.00%  "
# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Value = "'250.71"
$ws.Range("E5").Value = "  -2.68%  "
# Row 6
$ws.Range("D6").Value = "'0.7202"
$ws.Range("E6").Value = "  -7.64%  "
# Row 7
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  +0.07%  "
# Row 8
$ws.Range("D8").Value = "0.3367"
$ws.Range("E8").Value = "  -4.55%  "
# Row 9
$ws.Range("D9").Value = "'28.92"
$ws.Range("E9").Value = "  -0.67%  "
# Row 10
$ws.Range("D10").Value = "'0.07330"
$ws.Range("E10").Value = "  +4.12%  "
# Row 11
$ws.Range("D11").Value = "'0.8184"
$ws.Range("E11").Value = "  -5.65%  "
# Row 12
$ws.Range("D12").Value = "'0.08151"
$ws.Range("E12").Value = "  -0.77%  "
# Row 13
$ws.Range("D13").Value = "1.938.06"
$ws.Range("E13").Value = "  -3.01%  "
# Row 14
$ws.Range("D14").Value = "5.532"
$ws.Range("E14").Value = "  -1.09%  "
# Row 15
$ws.Range("D15").Value = "'95.31"
$ws.Range("E15").Value = "  -5.53%  "
# Row 16
$ws.Range("D16").Value = "'14.90"
$ws.Range("E16").Value = "  -3.96%  "
# Row 17
$ws.Range("D17").Value = "30.350.28"
$ws.Range("E17").Value = "  -2.76%  "
# Row 18
$ws.Range("D18").Value = "'0.000008275"
$ws.Range("E18").Value = "  +4.03%  "
# Row 19
$ws.Range("D19").Value = "'254.55"
$ws.Range("E19").Value = "  -7.22%  "
# Row 20
$ws.Range("D20").Value = "'5.890"
$ws.Range("E20").Value = "  -1.33%  "
# Row 21
$ws.Range("D21").Value = "2.193.94"
$ws.Range("E21").Value = "  -2.78%  "
# Row 22
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  +0.12%  "
# Row 23
$ws.Range("D23").Value = "'0.9997"
$ws.Range("E23").Value = "  -0.01%  "
# Row 24
$ws.Range("D24").Value = "'6.974"
$ws.Range("E24").Value = "  -2.50%  "
# Row 25
$ws.Range("D25").Value = "'9.862"
$ws.Range("E25").Value = "  -2.43%  "
# Row 26
$ws.Range("D26").Value = "160.37"
$ws.Range("E26").Value = "  -2.41%  "
# Row 27
$ws.Range("D27").Value = "'2.441"
$ws.Range("E27").Value = "  +3.10%  "
# Row 28
$ws.Range("D28").Value = "'19.46"
$ws.Range("E28").Value = "  -2.67%  "
# Row 29
$ws.Range("D29").Value = "0.1323"
$ws.Range("E29").Value = "  -11.14%  "
# Row 30
$ws.Range("E30").Value = "  -2.72%  "
# Row 31
$ws.Range("E31").Value = "  -0.66%  "
# Row 32
$ws.Range("D32").Value = "4.492"
$ws.Range("E32").Value = "  -2.83%  "
# Row 33
$ws.Range("D33").Value = "'4.248"
$ws.Range("E33").Value = "  -4.42%  "
# Row 34
$ws.Range("D34").Value = "'0.05249"
$ws.Range("E34").Value = "  +0.48%  "
# Row 35
$ws.Range("D35").Value = "1.272"
$ws.Range("E35").Value = "  +3.15%  "
# Row 36
$ws.Range("D36").Value = "'0.7548"
$ws.Range("E36").Value = "  -2.74%  "
# Row 37
$ws.Range("D37").Value = "2.738"
$ws.Range("E37").Value = "  -2.73%  "
# Row 38
$ws.Range("E38").Value = "  -1.02%  "
# Row 39
$ws.Range("E39").Value = "  -2.48%  "
# Row 40
$ws.Range("D40").Value = "'81.20"
$ws.Range("E40").Value = "  +1.70%  "
# Row 41
$ws.Range("D41").Value = "'6.580"
$ws.Range("E41").Value = "  -2.41%  "
# Row 42
$ws.Range("D42").Value = "'0.4583"
$ws.Range("E42").Value = "  -3.05%  "
# Row 43
$ws.Range("D43").Value = "'2.029"
$ws.Range("E43").Value = "  -5.72%  "
# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8439"
$ws.Range("E44").Value = "  -1.11%  "
# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  +0.06%  "
# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'102.49"
$ws.Range("E46").Value = "  -3.70%  "
# Row 47
$ws.Range("D47").Value = "'9.871"
$ws.Range("E47").Value = "  -1.48%  "
# Row 48
$ws.Range("D48").Value = "'7.449"
$ws.Range("E48").Value = "  -3.78%  "
# Row 49
$ws.Range("D49").Value = "'37.06"
$ws.Range("E49").Value = "  +0.51%  "
# Row 50
$ws.Range("D50").Value = "'0.4201"
$ws.Range("E50").Value = "  -3.23%  "
# Row 51
$ws.Range("D51").Value = "'1.511"
$ws.Range("E51").Value = "  -0.14%  "
